# Replace the literal password text with the placeholder instruction.
# The original paragraph is split across two runs:
#   "password: " + "xvaYs4rbkhHlMrtX"
# The target paragraph is a single run:
#   "Password: ask the trainer"
$d = $word.ActiveDocument

$d.Content.Find.Execute("password: xvaYs4rbkhHlMrtX", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Password: ask the trainer", 2)
